$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: reorder the block-name headers ---
# Old order: bedrooms_1, kitchens_1, bedrooms_2, kitchens_2, living_rooms_1, living_rooms_2
# New order: living_rooms_1, bedrooms_1, kitchens_1, living_rooms_2, bedrooms_2, kitchens_2
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# --- Rows 2-7: updated one-hot block-order values, re-aligned to the new columns ---
$data = @(
    @(0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 1, 0, 0, 0),
    @(1, 0, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}

$wb.Save()
